$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 46.00640566666667
$ws.Range("H2").Value = 138.019217
$ws.Range("I2").Value = 0.1268639696980132
$ws.Range("J2").Value = 0.1268639696980132
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.110028
$ws.Range("N2").Value = 0.330084
$ws.Range("Q2").Value = 5.061992802692
$ws.Range("R2").Value = 45.557935224228
$ws.Range("S2").Value = 0.1268639696980132
$ws.Range("T2").Value = 0.1268639696980132

# Row 3
$ws.Range("I3").Value = 0.3141149221428037
$ws.Range("J3").Value = 0.3141149221428037
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.110028
$ws.Range("N3").Value = 0.330084
$ws.Range("Q3").Value = 12.533483532716
$ws.Range("R3").Value = 112.801351794444
$ws.Range("S3").Value = 0.3141149221428037
$ws.Range("T3").Value = 0.3141149221428037

# Row 4
$ws.Range("G4").Value = 61.800369
$ws.Range("H4").Value = 185.401107
$ws.Range("I4").Value = 0.1704162719632449
$ws.Range("J4").Value = 0.1704162719632449
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.110028
$ws.Range("N4").Value = 0.330084
$ws.Range("Q4").Value = 6.799771000332
$ws.Range("R4").Value = 61.19793900298799
$ws.Range("S4").Value = 0.1704162719632449
$ws.Range("T4").Value = 0.1704162719632449

# Row 5
$ws.Range("G5").Value = 10.80043633333334
$ws.Range("H5").Value = 32.401309
$ws.Range("I5").Value = 0.02978250980189204
$ws.Range("J5").Value = 0.02978250980189203
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.110028
$ws.Range("N5").Value = 0.330084
$ws.Range("Q5").Value = 1.188350408884
$ws.Range("R5").Value = 10.695153679956
$ws.Range("S5").Value = 0.02978250980189204
$ws.Range("T5").Value = 0.02978250980189203

# Row 6
$ws.Range("G6").Value = 38.092607
$ws.Range("H6").Value = 114.277821
$ws.Range("I6").Value = 0.1050414452104809
$ws.Range("J6").Value = 0.1050414452104809
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.110028
$ws.Range("N6").Value = 0.330084
$ws.Range("Q6").Value = 4.191253362996
$ws.Range("R6").Value = 37.721280266964
$ws.Range("S6").Value = 0.1050414452104809
$ws.Range("T6").Value = 0.1050414452104809

# Row 7
$ws.Range("G7").Value = 92.03201033333333
$ws.Range("H7").Value = 276.096031
$ws.Range("I7").Value = 0.2537808811835653
$ws.Range("J7").Value = 0.2537808811835653
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.110028
$ws.Range("N7").Value = 0.330084
$ws.Range("Q7").Value = 10.126098032956
$ws.Range("R7").Value = 91.13488229660399
$ws.Range("S7").Value = 0.2537808811835653
$ws.Range("T7").Value = 0.2537808811835653
